$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the merge across A3:A6 so row 4 becomes independently editable
$ws.Range("A3:A6").UnMerge()

# Fill in row 4 with the new "80/50" epoch entry and its metrics
$ws.Range("A4").Value = "80/50"
$ws.Range("B4").Value = "transformer simple"

$ws.Range("C4").Value = 0.60591153580771195
$ws.Range("D4").Value = 0.76717203368462905
$ws.Range("E4").Value = 14.4484940182891
$ws.Range("F4").Value = 0.34216573947483497
$ws.Range("G4").Value = 0.52420581224005502
$ws.Range("H4").Value = 0.82352536484161898
$ws.Range("I4").Value = 16.197356524641201
$ws.Range("J4").Value = 0.44115773940263903

# Match the "transformer simple" row style (Consolas 7pt) used elsewhere on this sheet
$ws.Range("G3").Copy()
$ws.Range("C4:J4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update selection to reflect final cursor position
$ws.Range("C5").Select()
